# model-flow.pptx update
#  - bump the cached "datetimeFigureOut" date placeholder text on the
#    slide master + every slide layout from 7/7/2021 -> 7/14/2021
#  - on slide 1, swap the M2/M3 textbox labels (and their slightly
#    widened boxes) so the box that used to say "M2" now says "M3"
#    and vice versa

$p = $ppt.ActivePresentation

$oldDate = "7/7/2021"
$newDate = "7/14/2021"

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        if ($shape.PlaceholderFormat.Type -eq 16) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        Update-DateShape $layout.Shapes.Item($si)
    }
}

# Slide 1: swap the "M2"/"M3" textbox labels and bump their widths.
$slide1 = $p.Slides.Item(1)
$newWidthPt = 39.6586

for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $label = $shape.TextFrame.TextRange.Text
        if ($label -eq "M2") {
            $shapeM2 = $shape
        } elseif ($label -eq "M3") {
            $shapeM3 = $shape
        }
    }
}

if ($shapeM2 -and $shapeM3) {
    $shapeM2.TextFrame.TextRange.Text = "M3"
    $shapeM3.TextFrame.TextRange.Text = "M2"
    $shapeM2.Width = $newWidthPt
    $shapeM3.Width = $newWidthPt
}
